# Auto-generated cell value updates for Seraph_Profits workbook
# (scheduled market-price refresh - plain data values, no formulas)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 916.25
$ws.Range("I12").Value = 920
$ws.Range("J12").Value = 912.5
$ws.Range("K12").Value = 920
$ws.Range("L12").Value = 912.5
$ws.Range("M12").Value = -750
$ws.Range("N12").Value = -1252.5
$ws.Range("H40").Value = 2155.7144
$ws.Range("J40").Value = 2218
$ws.Range("L40").Value = 2218
$ws.Range("N40").Value = -2568
$ws.Range("H53").Value = 96.545456
$ws.Range("I53").Value = 131
$ws.Range("J53").Value = 55.2
$ws.Range("K53").Value = 131
$ws.Range("L53").Value = 55.2
$ws.Range("M53").Value = 506
$ws.Range("N53").Value = -1329.2
$ws.Range("H62").Value = 5796.4614
$ws.Range("J62").Value = 6785.7144
$ws.Range("L62").Value = 6785.7144
$ws.Range("N62").Value = -8033.7144
$ws.Range("H65").Value = 5796.4614
$ws.Range("J65").Value = 6785.7144
$ws.Range("L65").Value = 33928.572
$ws.Range("N65").Value = -40168.572
$ws.Range("H92").Value = 545.4400000000001
$ws.Range("I92").Value = 481.65
$ws.Range("J92").Value = 800.6
$ws.Range("K92").Value = 481.65
$ws.Range("L92").Value = 800.6
$ws.Range("M92").Value = 766.35
$ws.Range("N92").Value = -3296.6
$ws.Range("H98").Value = 847
$ws.Range("I98").Value = 650.7273
$ws.Range("K98").Value = 650.7273
$ws.Range("M98").Value = 847.2727
$ws.Range("H103").Value = 999.5
$ws.Range("J103").Value = 999.5
$ws.Range("L103").Value = 2998.5
$ws.Range("N103").Value = -4170.5
$ws.Range("H107").Value = 280.0909
$ws.Range("I107").Value = 280.0909
$ws.Range("K107").Value = 280.0909
$ws.Range("M107").Value = 1639.9091
$ws.Range("H111").Value = 2632.2
$ws.Range("H122").Value = 847
$ws.Range("I122").Value = 650.7273
$ws.Range("K122").Value = 1952.1819
$ws.Range("M122").Value = 497.8181
$ws.Range("H132").Value = 2332.3333
$ws.Range("I132").Value = 2332.3333
$ws.Range("K132").Value = 6996.999899999999
$ws.Range("M132").Value = -4466.999899999999
$ws.Range("H135").Value = 2717.2856
$ws.Range("I135").Value = 2099.5
$ws.Range("J135").Value = 2964.4
$ws.Range("K135").Value = 18895.5
$ws.Range("L135").Value = 26679.6
$ws.Range("M135").Value = -16360.5
$ws.Range("N135").Value = -31749.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4685.364
$ws.Range("I61").Value = 4616.2856
$ws.Range("K61").Value = 4616.2856
$ws.Range("M61").Value = -4404.2856
$ws.Range("H102").Value = 1170.8572
$ws.Range("I102").Value = 949.3333
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 949.3333
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = 672.6667
$ws.Range("N102").Value = -5744
$ws.Range("H136").Value = 4685.364
$ws.Range("I136").Value = 4616.2856
$ws.Range("K136").Value = 13848.8568
$ws.Range("M136").Value = -11298.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4536.9
$ws.Range("I105").Value = 3461.7058
$ws.Range("K105").Value = 3461.7058
$ws.Range("M105").Value = -1714.7058
$ws.Range("H107").Value = 618.13794
$ws.Range("I107").Value = 608.7308
$ws.Range("K107").Value = 608.7308
$ws.Range("M107").Value = 1311.2692

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = ""
$ws.Range("H16").Value = 3294.7368
$ws.Range("I16").Value = 3116.875
$ws.Range("J16").Value = 4243.3335
$ws.Range("K16").Value = 3116.875
$ws.Range("L16").Value = 4243.3335
$ws.Range("M16").Value = -2829.875
$ws.Range("N16").Value = -4817.3335
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = ""
$ws.Range("H99").Value = 13707
$ws.Range("I99").Value = 11040.833
$ws.Range("K99").Value = 11040.833
$ws.Range("M99").Value = -9542.833000000001
$ws.Range("H113").Value = 3294.7368
$ws.Range("I113").Value = 3116.875
$ws.Range("J113").Value = 4243.3335
$ws.Range("K113").Value = 3116.875
$ws.Range("L113").Value = 4243.3335
$ws.Range("M113").Value = -946.875
$ws.Range("N113").Value = -8583.333500000001
$ws.Range("H126").Value = 13707
$ws.Range("I126").Value = 11040.833
$ws.Range("K126").Value = 33122.499
$ws.Range("M126").Value = -30652.499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 13725
$ws.Range("I32").Value = 1462.5
$ws.Range("J32").Value = 25987.5
$ws.Range("K32").Value = 4387.5
$ws.Range("L32").Value = 77962.5
$ws.Range("M32").Value = -4104.5
$ws.Range("N32").Value = -78528.5
$ws.Range("H46").Value = 3335100
$ws.Range("I46").Value = 300
$ws.Range("K46").Value = 900
$ws.Range("M46").Value = -809
$ws.Range("H55").Value = 127049.25
$ws.Range("I55").Value = 250598.75
$ws.Range("J55").Value = 3499.75
$ws.Range("K55").Value = 751796.25
$ws.Range("L55").Value = 10499.25
$ws.Range("M55").Value = -751619.25
$ws.Range("N55").Value = -10853.25
$ws.Range("H114").Value = 1819.5
$ws.Range("I114").Value = 623.5
$ws.Range("J114").Value = 3015.5
$ws.Range("K114").Value = 1870.5
$ws.Range("L114").Value = 9046.5
$ws.Range("M114").Value = 1383.5
$ws.Range("N114").Value = -15554.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = ""
$ws.Range("H70").Value = 4444.1113
$ws.Range("I70").Value = 2498.5
$ws.Range("K70").Value = 2498.5
$ws.Range("M70").Value = -2228.5
$ws.Range("H73").Value = 4444.1113
$ws.Range("I73").Value = 2498.5
$ws.Range("K73").Value = 2498.5
$ws.Range("M73").Value = -1562.5
$ws.Range("H107").Value = 1005.913
$ws.Range("I107").Value = 1110.7
$ws.Range("K107").Value = 1110.7
$ws.Range("M107").Value = 809.3
$ws.Range("H113").Value = 5791.091
$ws.Range("I113").Value = 7263
$ws.Range("K113").Value = 7263
$ws.Range("M113").Value = -5093
$ws.Range("H122").Value = 74706.14
$ws.Range("I122").Value = 2476.889
$ws.Range("J122").Value = 204718.8
$ws.Range("K122").Value = 7430.667
$ws.Range("L122").Value = 614156.3999999999
$ws.Range("M122").Value = -4980.667
$ws.Range("N122").Value = -619056.3999999999
$ws.Range("H126").Value = 4582.8184
$ws.Range("I126").Value = 3899
$ws.Range("K126").Value = 11697
$ws.Range("M126").Value = -9227
$ws.Range("H132").Value = 2847.3
$ws.Range("I132").Value = 1337.25
$ws.Range("K132").Value = 4011.75
$ws.Range("M132").Value = -1481.75
$ws.Range("H136").Value = 24052.666
$ws.Range("J136").Value = 24052.666
$ws.Range("L136").Value = 72157.99800000001
$ws.Range("N136").Value = -77257.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2825
$ws.Range("I7").Value = 2825
$ws.Range("K7").Value = 2825
$ws.Range("M7").Value = -2713
$ws.Range("H40").Value = 3804.6667
$ws.Range("I40").Value = 3704
$ws.Range("K40").Value = 3704
$ws.Range("M40").Value = -3568
$ws.Range("H55").Value = 362.06668
$ws.Range("J55").Value = 286.66666
$ws.Range("L55").Value = 286.66666
$ws.Range("N55").Value = -632.66666
$ws.Range("H93").Value = 1996.5625
$ws.Range("I93").Value = 1996.5625
$ws.Range("K93").Value = 1996.5625
$ws.Range("M93").Value = -748.5625
$ws.Range("H126").Value = 2825
$ws.Range("I126").Value = 2825
$ws.Range("K126").Value = 8475
$ws.Range("M126").Value = -6005
$ws.Range("H132").Value = 1996
$ws.Range("I132").Value = 1996.3334
$ws.Range("K132").Value = 5989.0002
$ws.Range("M132").Value = -3459.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 35676
$ws.Range("J40").Value = 41014
$ws.Range("L40").Value = 41014
$ws.Range("N40").Value = -41312
$ws.Range("H126").Value = 1507.4
$ws.Range("I126").Value = 1049.5883
$ws.Range("J126").Value = 4101.6665
$ws.Range("K126").Value = 3148.7649
$ws.Range("L126").Value = 12304.9995
$ws.Range("M126").Value = -678.7648999999997
$ws.Range("N126").Value = -17244.9995

Write-Output "Updated Seraph_Profits sheets"
